$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.007.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.298.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.43%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.290.16'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.36%  '

$ws.Range("E8").Value = '  -4.91%  '

$ws.Range("E10").Value = '  -1.48%  '

$ws.Range("E11").Value = '  +15.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '38.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("E13").Value = '  -0.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.829.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.307.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.737.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.980'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("E21").Value = '  +1.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '294.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.173'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("E31").Value = '  -2.52%  '

$ws.Range("E32").Value = '  -0.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.110'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.82%  '

$ws.Range("E34").Value = '  -2.90%  '

$ws.Range("E35").Value = '  +14.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0477'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '134.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("E43").Value = '  -2.56%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.46%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.10%  '

$ws.Range("E48").Value = '  +2.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.103.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.635.04'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.85%  '
